# Append the new listings scraped from slando (rows 56-110) to the "Rent Data" sheet.
# Row 56 only refreshes the date (col E) on a previously-seen aviso listing; rows 86-110
# are brand-new slando listings. Column layout: A=source, B=title/location, C=description/title,
# D=price, E=date, F=url, G=constant marker.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row 56
$ws.Cells.Item(56,1).Value = 'aviso'
$ws.Cells.Item(56,2).Value = 'Сдам уютную квартиру на Академика Вильямса 81/2'
$ws.Cells.Item(56,3).Value = ' Квартира в хорошем жилом состоянии.Все необходимое для жизни. Без балкона. Свое отопление, бойлер. 2200 грн.'
$ws.Cells.Item(56,4).Value = ' 2 200 грн. в месяц'
$ws.Cells.Item(56,5).Value = '21-12-2013 / 13:49:52'
$ws.Cells.Item(56,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5371361'
$ws.Cells.Item(56,7).Value = 'test field other'

# row 57
$ws.Cells.Item(57,1).Value = 'aviso'
$ws.Cells.Item(57,2).Value = 'Сдам 1 ком. квартиру на Черемушках »'
$ws.Cells.Item(57,3).Value = ' Сдам 1 ком. квартиру на Черемушках, ул. Гайдара/Вишневского. 3/5 эт. дома, 32/18/6. МПО, бронированная дверь, ремонт, балкон застеклен. В комнате и кухне хорошая мебель и вся бытовая техника, бойлер...'
$ws.Cells.Item(57,4).Value = ' 2 700 грн. в месяц'
$ws.Cells.Item(57,5).Value = '21-12-2013 / 10:55:11'
$ws.Cells.Item(57,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696444'
$ws.Cells.Item(57,7).Value = 'test field other'

# row 58
$ws.Cells.Item(58,1).Value = 'aviso'
$ws.Cells.Item(58,2).Value = 'Сдам в длительную аренду 2 к »'
$ws.Cells.Item(58,3).Value = ' Сдам в длительную аренду 2 к квартиру на Вузовском. 2 -ком. квартира со всеми удобствами на Вузовском. Ремонт 2013 года, мет. пластиковое остекление, все новае: мебель, двухспальная кровать, диван, встроенная...'
$ws.Cells.Item(58,4).Value = ' 500 у.е. в месяц'
$ws.Cells.Item(58,5).Value = '21-12-2013 / 10:55:11'
$ws.Cells.Item(58,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696133'
$ws.Cells.Item(58,7).Value = 'test field other'

# row 59
$ws.Cells.Item(59,1).Value = 'aviso'
$ws.Cells.Item(59,2).Value = 'Сдам кв-ру Королева/Левитана »'
$ws.Cells.Item(59,3).Value = ' Сдам 2-х комнатную квартиру с евроремонтом, мебелью и техникой - Королева/Левитана 2/5, общая площадь 45 кв.м, Жилая 36.6 кв.м. Высота потолков 2.7м. Ремонту 3 года, интерьер выполнен в светлых тонах. На полу...'
$ws.Cells.Item(59,4).Value = ' 400 у.е. в месяц'
$ws.Cells.Item(59,5).Value = '21-12-2013 / 09:33:50'
$ws.Cells.Item(59,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5722413'
$ws.Cells.Item(59,7).Value = 'test field other'

# row 60
$ws.Cells.Item(60,1).Value = 'aviso'
$ws.Cells.Item(60,2).Value = 'квартира Центр, жк"Капитан"/Парк »'
$ws.Cells.Item(60,3).Value = ' 9/20 эт. 95/50/23. Стильная квартира после качественного ремонта. 2 отдельные спальни и кухня-студия. Встроенная мебель. Бытовая техник Bosch - посудомоечная и стиральная машина, холодильник, духовой шкаф, 3...'
$ws.Cells.Item(60,4).Value = ' 800 у.е. в месяц'
$ws.Cells.Item(60,5).Value = '21-12-2013 / 09:03:44'
$ws.Cells.Item(60,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5722383'
$ws.Cells.Item(60,7).Value = 'test field other'

# row 61
$ws.Cells.Item(61,1).Value = 'aviso'
$ws.Cells.Item(61,2).Value = 'Бунина'
$ws.Cells.Item(61,3).Value = ' Сдам долгосрочно 2-к квартиру Одесса, Приморский, 4200 грн./месяц 2/3 Площадь общая: 52.00 жилая: 33.00 кухни: 7.00 кв.м. Приморский р-н.Бунина 39, 52/33/7, 2/3 эт., высота потолков 2,8м., квартира частая...'
$ws.Cells.Item(61,4).Value = ' 4 200 грн.'
$ws.Cells.Item(61,5).Value = '20-12-2013 / 19:17:39'
$ws.Cells.Item(61,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5722385'
$ws.Cells.Item(61,7).Value = 'test field other'

# row 62
$ws.Cells.Item(62,1).Value = 'aviso'
$ws.Cells.Item(62,2).Value = 'Сдам 2-комн., Королева/ р-к «Южный»'
$ws.Cells.Item(62,3).Value = ' капитальный ремонт, МПО, современная облицовка, бойлер, кондиционер, машина-автомат, холодильник, телевизор, встроенная кухня, шкафы-купе, вся необходимая мебель, 3 000 грн./мес'
$ws.Cells.Item(62,4).Value = ' 3 000 грн. в месяц'
$ws.Cells.Item(62,5).Value = '20-12-2013 / 16:10:24'
$ws.Cells.Item(62,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696670'
$ws.Cells.Item(62,7).Value = 'test field other'

# row 63
$ws.Cells.Item(63,1).Value = 'aviso'
$ws.Cells.Item(63,2).Value = 'Котовского'
$ws.Cells.Item(63,3).Value = ' Сахарова ул., 5/11-эт. дома, 70 кв.м, элитный новый дом, стильный евроремонт 2013 года, новая итальянская мебель и брендовая элитная техника. 2900грн/мес. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(63,4).Value = ' 2 900 грн.'
$ws.Cells.Item(63,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(63,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696526'
$ws.Cells.Item(63,7).Value = 'test field other'

# row 64
$ws.Cells.Item(64,1).Value = 'aviso'
$ws.Cells.Item(64,2).Value = 'Аркадия'
$ws.Cells.Item(64,3).Value = ' Генуэзская ул., спальня, гостинная 26 м, шкаф-купе, новая мягкая мебель, акриловая ванна, панорамное остекление, кухня Skavollini, новая техника, 4100грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(64,4).Value = ' 4 100 грн.'
$ws.Cells.Item(64,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(64,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696525'
$ws.Cells.Item(64,7).Value = 'test field other'

# row 65
$ws.Cells.Item(65,1).Value = 'aviso'
$ws.Cells.Item(65,2).Value = 'Екатерининская ул.'
$ws.Cells.Item(65,3).Value = ' 8/9-эт. дома, кап.ремонт, р/х, 3 лоджии застеклены-МПЛ, кухня 9м, мягкий уголок, 2 телевизора, 2 спальная кровать, шкаф-купе, быт. техника, 3700грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(65,4).Value = ' 3 700 грн.'
$ws.Cells.Item(65,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(65,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696524'
$ws.Cells.Item(65,7).Value = 'test field other'

# row 66
$ws.Cells.Item(66,1).Value = 'aviso'
$ws.Cells.Item(66,2).Value = '6 ст.Б.Фонтана'
$ws.Cells.Item(66,3).Value = ' спальня и кухня-студия, дорогой авторский ремонт, брендовая быттехника, евро санузел, плазменный телевизор, гардеробная, ламинат, теплые полы, паркинг, 4200грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(66,4).Value = ' 4 200 грн.'
$ws.Cells.Item(66,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(66,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696522'
$ws.Cells.Item(66,7).Value = 'test field other'

# row 67
$ws.Cells.Item(67,1).Value = 'aviso'
$ws.Cells.Item(67,2).Value = 'Бунина ул./Ришельевская ул.'
$ws.Cells.Item(67,3).Value = ' 3/5-эт. дома, "бельгийка", евроремонт, р/х, итальянская спальня, кожанная мебель, 2 балкона застеклены, плазма, техника "Электролюкс", 4700грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(67,4).Value = ' 4 700 грн.'
$ws.Cells.Item(67,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(67,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696523'
$ws.Cells.Item(67,7).Value = 'test field other'

# row 68
$ws.Cells.Item(68,1).Value = 'aviso'
$ws.Cells.Item(68,2).Value = 'Гагарина пр./Шевченко пр.'
$ws.Cells.Item(68,3).Value = ' евроремонт 2013г., дизайнерская мебель, 2 с.у., джакузи, панорамное остекление, 2 плазменных телевизора, стильный интерьер, 4800грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(68,4).Value = ' 4 800 грн.'
$ws.Cells.Item(68,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(68,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696521'
$ws.Cells.Item(68,7).Value = 'test field other'

# row 69
$ws.Cells.Item(69,1).Value = 'aviso'
$ws.Cells.Item(69,2).Value = 'Шевченко пр.'
$ws.Cells.Item(69,3).Value = ' "Легион", авторский дизайн, шкаф-купе, итальянский 2 спальный гарнитур, брендовая быттехника, панорамное остекление, евро санузел, паркинг, охрана, 3600грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(69,4).Value = ' 3 600 грн.'
$ws.Cells.Item(69,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(69,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696520'
$ws.Cells.Item(69,7).Value = 'test field other'

# row 70
$ws.Cells.Item(70,1).Value = 'aviso'
$ws.Cells.Item(70,2).Value = 'Французский б-р. Каркашадзе'
$ws.Cells.Item(70,3).Value = ' евроремонт 2013, спальня+студия, итальянская мебель. панорамный вид, плазма, супер-стиль, гардеробная, посуточно, варианты от 2800, 4800грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(70,4).Value = ' 4 800 грн.'
$ws.Cells.Item(70,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(70,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696519'
$ws.Cells.Item(70,7).Value = 'test field other'

# row 71
$ws.Cells.Item(71,1).Value = 'aviso'
$ws.Cells.Item(71,2).Value = 'Гагаринское плато'
$ws.Cells.Item(71,3).Value = ' "Жемчужина", евроремонт, спальня и студия, дизайнерская мебель, 2 с.у., акриловая ванна, гардеробная, 2 плазменных телевизора, панорамное остекление, 4450грн/мес., торг. Тел: (50) 3912828, (67) 4862288'
$ws.Cells.Item(71,4).Value = ' 4 450 грн.'
$ws.Cells.Item(71,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(71,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696518'
$ws.Cells.Item(71,7).Value = 'test field other'

# row 72
$ws.Cells.Item(72,1).Value = 'aviso'
$ws.Cells.Item(72,2).Value = 'Глушко пр.'
$ws.Cells.Item(72,3).Value = ' 8/9, "чешка", после ремонта, спальный гарнитур, бытовая техника, кондиционер., 2800грн/мес. Тел: (48) 7949608, (97) 2575077'
$ws.Cells.Item(72,4).Value = ' 2 800 грн.'
$ws.Cells.Item(72,5).Value = '20-12-2013 / 15:59:28'
$ws.Cells.Item(72,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696308'
$ws.Cells.Item(72,7).Value = 'test field other'

# row 73
$ws.Cells.Item(73,1).Value = 'aviso'
$ws.Cells.Item(73,2).Value = 'Сдам 2-комн., Люстдорфская дорога/ »'
$ws.Cells.Item(73,3).Value = ' 6/9, 50/30/8, капитальный ремонт, МПО, предметы интерьера, встроенная кухня, бойлер 2 800 грн./ мес.'
$ws.Cells.Item(73,4).Value = ' 2 800 грн. в месяц'
$ws.Cells.Item(73,5).Value = '20-12-2013 / 14:24:40'
$ws.Cells.Item(73,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5612275'
$ws.Cells.Item(73,7).Value = 'test field other'

# row 74
$ws.Cells.Item(74,1).Value = 'aviso'
$ws.Cells.Item(74,2).Value = 'Сдам 2-комн.,Княжеская/ «Новый рынок»'
$ws.Cells.Item(74,3).Value = ' Сдам 2-х комн., квартиру на Княжеской/ «Новый рынок» евроремонт, 2-е раздельные комнаты, АГВ, вся необходимая мебель и бытовая техника, есть свой дворик, место для машины, 3 500 грн./мес. СРОЧНО'
$ws.Cells.Item(74,4).Value = ' 3 500 грн. в месяц'
$ws.Cells.Item(74,5).Value = '20-12-2013 / 14:24:27'
$ws.Cells.Item(74,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5695737'
$ws.Cells.Item(74,7).Value = 'test field other'

# row 75
$ws.Cells.Item(75,1).Value = 'aviso'
$ws.Cells.Item(75,2).Value = 'Екатерининская'
$ws.Cells.Item(75,3).Value = ' Сдам долгосрочно 1-к квартиру Одесса, Приморский, 450 $/месяц 4/4 Площадь общая: 34.00 жилая: 19.00 кухни: 9.00 кв.м. 1 комн. на Екатерининской/Базарной, дому 10 лет, красивый евроремонт, новая мебель...'
$ws.Cells.Item(75,4).Value = ' 450 у.е.'
$ws.Cells.Item(75,5).Value = '20-12-2013 / 13:57:11'
$ws.Cells.Item(75,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5695570'
$ws.Cells.Item(75,7).Value = 'test field other'

# row 76
$ws.Cells.Item(76,1).Value = 'aviso'
$ws.Cells.Item(76,2).Value = 'Прохоровская ул./Банный'
$ws.Cells.Item(76,3).Value = ' пер., 1/2, комн. 22 кв.м, кухня 7 кв.м, т/душ, АОГВ, стиральная машина-автомат, мебель, хол., длительно, 2400 грн. Тел: (48) 7043990, (67) 2842428, моб.'
$ws.Cells.Item(76,4).Value = ' 2 400 грн.'
$ws.Cells.Item(76,5).Value = '19-12-2013 / 18:00:00'
$ws.Cells.Item(76,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696470'
$ws.Cells.Item(76,7).Value = 'test field other'

# row 77
$ws.Cells.Item(77,1).Value = 'aviso'
$ws.Cells.Item(77,2).Value = 'Центр'
$ws.Cells.Item(77,3).Value = ' ул Пишоновская новострой ,сдам свою однокомнатную квартиру,евроремонт с новой мебелью ,стиралка ,кондиционер-бытовая техника.Раньше не сдавалась.цена 3200.тел 0635372827б.0672811079'
$ws.Cells.Item(77,4).Value = ' 3 200 грн.'
$ws.Cells.Item(77,5).Value = '19-12-2013 / 17:55:01'
$ws.Cells.Item(77,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696468'
$ws.Cells.Item(77,7).Value = 'test field other'

# row 78
$ws.Cells.Item(78,1).Value = 'aviso'
$ws.Cells.Item(78,2).Value = 'Сдам 2-х комнатную впервые новострой »'
$ws.Cells.Item(78,3).Value = ' Сдам 2-х комнатную впервые новострой ЖК"Тирас" на Таирово 25 Чапаевской Дивизии / А.Королева, 2 / 15, соврем.ремонт,общая площадь 75 кв.м.,комнаты с раздельными ходами 20 и 18 кв.м.,кухня 14 кв.м..с / т...'
$ws.Cells.Item(78,4).Value = ' 600 у.е. в месяц'
$ws.Cells.Item(78,5).Value = '19-12-2013 / 16:17:27'
$ws.Cells.Item(78,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696443'
$ws.Cells.Item(78,7).Value = 'test field other'

# row 79
$ws.Cells.Item(79,1).Value = 'aviso'
$ws.Cells.Item(79,2).Value = 'Молдаванка, Прохоровская'
$ws.Cells.Item(79,3).Value = ' ул., 2/4, кирпич, ремонт, комната 25 кв.м, кухня 10 кв.м, мебель, бытовая техника, возможно без мебели. Варианты от 2200 грн. Тел: (97) 2952343, (63) 6403057'
$ws.Cells.Item(79,4).Value = ' 2 200 грн.'
$ws.Cells.Item(79,5).Value = '19-12-2013 / 16:00:00'
$ws.Cells.Item(79,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696434'
$ws.Cells.Item(79,7).Value = 'test field other'

# row 80
$ws.Cells.Item(80,1).Value = 'aviso'
$ws.Cells.Item(80,2).Value = 'Новикова ул., Застава-2'
$ws.Cells.Item(80,3).Value = ' 2/5, бойлер, холодильник, стиральная машина. 2000 грн. Варианты. Тел: (48) 7948963, (67) 9278963, (93) 8588821'
$ws.Cells.Item(80,4).Value = ' 2 000 грн.'
$ws.Cells.Item(80,5).Value = '19-12-2013 / 16:00:00'
$ws.Cells.Item(80,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696436'
$ws.Cells.Item(80,7).Value = 'test field other'

# row 81
$ws.Cells.Item(81,1).Value = 'aviso'
$ws.Cells.Item(81,2).Value = 'Базарная/Белинского ул.'
$ws.Cells.Item(81,3).Value = ' 2/2, комн. разд., кухня 10 кв.м, свежий капремонт, с мебелью и бытовой техникой, т/в, двухконтурн. котел. Сдается впервые. 3200 грн. Тел: (48) 7946045, (67) 4885491'
$ws.Cells.Item(81,4).Value = ' 3 200 грн.'
$ws.Cells.Item(81,5).Value = '19-12-2013 / 16:00:00'
$ws.Cells.Item(81,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696426'
$ws.Cells.Item(81,7).Value = 'test field other'

# row 82
$ws.Cells.Item(82,1).Value = 'aviso'
$ws.Cells.Item(82,2).Value = 'Шевченко пр./Шампанский'
$ws.Cells.Item(82,3).Value = ' пер., евроремонт, новая мебель и бытовая техника, Wi-Fi, теплая, светлая и уютная квартира, хозяин, 3300 грн. Тел: (93) 0346098, моб.'
$ws.Cells.Item(82,4).Value = ' 3 300 грн.'
$ws.Cells.Item(82,5).Value = '19-12-2013 / 15:00:00'
$ws.Cells.Item(82,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696411'
$ws.Cells.Item(82,7).Value = 'test field other'

# row 83
$ws.Cells.Item(83,1).Value = 'aviso'
$ws.Cells.Item(83,2).Value = 'Большая Арнаутская'
$ws.Cells.Item(83,3).Value = ' ул./Александровский пр., 76/32+18/26, h=4м, двухконт. котел, плазма, вся бытовая техника, итальянский раскладывающийся уголок, двухспальная кровать, шкафы-купе, встроенная кухня, Интернет, Wi-Fi, каб.ТВ...'
$ws.Cells.Item(83,4).Value = ' 5 600 грн.'
$ws.Cells.Item(83,5).Value = '19-12-2013 / 15:00:00'
$ws.Cells.Item(83,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696395'
$ws.Cells.Item(83,7).Value = 'test field other'

# row 84
$ws.Cells.Item(84,1).Value = 'aviso'
$ws.Cells.Item(84,2).Value = 'Вильямса ул., Бастма'
$ws.Cells.Item(84,3).Value = ' ", 5/11, комн. разд., новострой, евроремонт 2013г., новая стильная мебель и вся бытовая техника, кухня-студия + 2 комнаты. Сдается впервые. Варианты. 5000 грн. Срочно! Тел: (67) 4853097'
$ws.Cells.Item(84,4).Value = ' 5 000 грн.'
$ws.Cells.Item(84,5).Value = '19-12-2013 / 15:00:00'
$ws.Cells.Item(84,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696401'
$ws.Cells.Item(84,7).Value = 'test field other'

# row 85
$ws.Cells.Item(85,1).Value = 'aviso'
$ws.Cells.Item(85,2).Value = 'Адмиральский пр., 5 ст.Б.'
$ws.Cells.Item(85,3).Value = ' Фонтана, 3/5, евроремонт 2013г., кухня- студия 24 кв.м + спальня, лоджия 14 кв.м, красивая, стильная, молодежная, все новое. 4500 грн., торг. Срочно! Тел: (67) 4853097'
$ws.Cells.Item(85,4).Value = ' 4 500 грн.'
$ws.Cells.Item(85,5).Value = '19-12-2013 / 15:00:00'
$ws.Cells.Item(85,6).Value = 'http://www.aviso.ua/odessa/view.php?adid=5696402'
$ws.Cells.Item(85,7).Value = 'test field other'

# row 86
$ws.Cells.Item(86,1).Value = 'slando'
$ws.Cells.Item(86,2).Value = ', Приморский'
$ws.Cells.Item(86,3).Value = '1 комн. в центре на ул. Приморская(008)'
$ws.Cells.Item(86,4).Value = '3 500 грн.'
$ws.Cells.Item(86,5).Value = 'Сегодня 13:51'
$ws.Cells.Item(86,6).Value = 'http://odessa.od.slando.ua/obyavlenie/1-komn-v-tsentre-na-ul-primorskaya008-ID87HOw.html'
$ws.Cells.Item(86,7).Value = 'test field other'

# row 87
$ws.Cells.Item(87,1).Value = 'slando'
$ws.Cells.Item(87,2).Value = ', Суворовский, Посёлок Котовского'
$ws.Cells.Item(87,3).Value = '2-х квартира на Днепродороге парковая зона'
$ws.Cells.Item(87,4).Value = '3 100 грн.'
$ws.Cells.Item(87,5).Value = 'Сегодня 13:49'
$ws.Cells.Item(87,6).Value = 'http://odessa.od.slando.ua/obyavlenie/2-h-kvartira-na-dneprodoroge-parkovaya-zona-ID7PWSz.html'
$ws.Cells.Item(87,7).Value = 'test field other'

# row 88
$ws.Cells.Item(88,1).Value = 'slando'
$ws.Cells.Item(88,2).Value = ', Приморский, Ольгиевская / Княжеская'
$ws.Cells.Item(88,3).Value = 'Сдам 1 комнатную квартиру'
$ws.Cells.Item(88,4).Value = '2 500 грн.'
$ws.Cells.Item(88,5).Value = 'Сегодня 13:47'
$ws.Cells.Item(88,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-1-komnatnuyu-kvartiru-ID7oScN.html'
$ws.Cells.Item(88,7).Value = 'test field other'

# row 89
$ws.Cells.Item(89,1).Value = 'slando'
$ws.Cells.Item(89,2).Value = ', Академика Филатова'
$ws.Cells.Item(89,3).Value = 'сдам 2 квартиру Филатова/солнечное'
$ws.Cells.Item(89,4).Value = '2 800 грн.'
$ws.Cells.Item(89,5).Value = 'Сегодня 13:46'
$ws.Cells.Item(89,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-kvartiru-filatova-solnechnoe-ID5jafj.html'
$ws.Cells.Item(89,7).Value = 'test field other'

# row 90
$ws.Cells.Item(90,1).Value = 'slando'
$ws.Cells.Item(90,2).Value = ', Приморский, Тираспольская'
$ws.Cells.Item(90,3).Value = 'сда 2х.к.кв тираспольская'
$ws.Cells.Item(90,4).Value = '4 000 грн.'
$ws.Cells.Item(90,5).Value = 'Сегодня 13:46'
$ws.Cells.Item(90,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sda-2h-k-kv-tiraspolskaya-ID6Uglr.html'
$ws.Cells.Item(90,7).Value = 'test field other'

# row 91
$ws.Cells.Item(91,1).Value = 'slando'
$ws.Cells.Item(91,2).Value = ', Приморский'
$ws.Cells.Item(91,3).Value = '2х комн. кв. в центре Ришельевская (004)'
$ws.Cells.Item(91,4).Value = '6 500 грн.'
$ws.Cells.Item(91,5).Value = 'Сегодня 13:44'
$ws.Cells.Item(91,6).Value = 'http://odessa.od.slando.ua/obyavlenie/2h-komn-kv-v-tsentre-rishelevskaya-004-ID87Ab3.html'
$ws.Cells.Item(91,7).Value = 'test field other'

# row 92
$ws.Cells.Item(92,1).Value = 'slando'
$ws.Cells.Item(92,2).Value = ', Еврейская/Осипова'
$ws.Cells.Item(92,3).Value = 'Сдам квартиру-студию в центре'
$ws.Cells.Item(92,4).Value = '3 200 грн.'
$ws.Cells.Item(92,5).Value = 'Сегодня 13:44'
$ws.Cells.Item(92,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-kvartiru-studiyu-v-tsentre-ID5AuP9.html'
$ws.Cells.Item(92,7).Value = 'test field other'

# row 93
$ws.Cells.Item(93,1).Value = 'slando'
$ws.Cells.Item(93,2).Value = ', Приморский, греческая'
$ws.Cells.Item(93,3).Value = 'Сдам 2-Х ком.кв.центр'
$ws.Cells.Item(93,4).Value = '6 511 грн.'
$ws.Cells.Item(93,5).Value = 'Сегодня 13:42'
$ws.Cells.Item(93,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-h-kom-kv-tsentr-ID8b2IZ.html'
$ws.Cells.Item(93,7).Value = 'test field other'

# row 94
$ws.Cells.Item(94,1).Value = 'slando'
$ws.Cells.Item(94,2).Value = ', Киевский, Королёва/Костанди'
$ws.Cells.Item(94,3).Value = 'Сдам квартиру кухня-студия+спальня на Королёва'
$ws.Cells.Item(94,4).Value = '4 069 грн.'
$ws.Cells.Item(94,5).Value = 'Сегодня 13:42'
$ws.Cells.Item(94,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-kvartiru-kuhnya-studiya-spalnya-na-koroleva-ID87IST.html'
$ws.Cells.Item(94,7).Value = 'test field other'

# row 95
$ws.Cells.Item(95,1).Value = 'slando'
$ws.Cells.Item(95,2).Value = ', Приморский'
$ws.Cells.Item(95,3).Value = 'Сдам квартиру-студию на Раскидайловской'
$ws.Cells.Item(95,4).Value = '2 500 грн.'
$ws.Cells.Item(95,5).Value = 'Сегодня 13:42'
$ws.Cells.Item(95,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-kvartiru-studiyu-na-raskidaylovskoy-ID72zJX.html'
$ws.Cells.Item(95,7).Value = 'test field other'

# row 96
$ws.Cells.Item(96,1).Value = 'slando'
$ws.Cells.Item(96,2).Value = ', Приморский, Б Арнаутская'
$ws.Cells.Item(96,3).Value = 'Сдам 2-Х ком.кв.'
$ws.Cells.Item(96,4).Value = '4 500 грн.'
$ws.Cells.Item(96,5).Value = 'Сегодня 13:41'
$ws.Cells.Item(96,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-h-kom-kv-ID7IFxd.html'
$ws.Cells.Item(96,7).Value = 'test field other'

# row 97
$ws.Cells.Item(97,1).Value = 'slando'
$ws.Cells.Item(97,2).Value = ', Приморский, среднефонтанская'
$ws.Cells.Item(97,3).Value = 'сдам 1 квартиру Среднефонтанская'
$ws.Cells.Item(97,4).Value = '2 900 грн.'
$ws.Cells.Item(97,5).Value = 'Сегодня 13:40'
$ws.Cells.Item(97,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-1-kvartiru-srednefontanskaya-ID6Oez3.html'
$ws.Cells.Item(97,7).Value = 'test field other'

# row 98
$ws.Cells.Item(98,1).Value = 'slando'
$ws.Cells.Item(98,2).Value = ', Приморский, Троицкая/Александровский пр-т'
$ws.Cells.Item(98,3).Value = 'Сдам 1-ную квартиру'
$ws.Cells.Item(98,4).Value = '6 104 грн.'
$ws.Cells.Item(98,5).Value = 'Сегодня 13:39'
$ws.Cells.Item(98,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-1-nuyu-kvartiru-ID8b3pW.html'
$ws.Cells.Item(98,7).Value = 'test field other'

# row 99
$ws.Cells.Item(99,1).Value = 'slando'
$ws.Cells.Item(99,2).Value = ', Малиновский, Варненская'
$ws.Cells.Item(99,3).Value = 'Сдам 2-х комнатную квартиру на Черёмушках'
$ws.Cells.Item(99,4).Value = '3 500 грн.'
$ws.Cells.Item(99,5).Value = 'Сегодня 13:36'
$ws.Cells.Item(99,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-h-komnatnuyu-kvartiru-na-cheremushkah-ID87IKz.html'
$ws.Cells.Item(99,7).Value = 'test field other'

# row 100
$ws.Cells.Item(100,1).Value = 'slando'
$ws.Cells.Item(100,2).Value = ', Суворовский, Посёлок Котовского'
$ws.Cells.Item(100,3).Value = 'хорошая 1 комн недорого Бочарова'
$ws.Cells.Item(100,4).Value = '1 800 грн.'
$ws.Cells.Item(100,5).Value = 'Сегодня 13:36'
$ws.Cells.Item(100,6).Value = 'http://odessa.od.slando.ua/obyavlenie/horoshaya-1-komn-nedorogo-bocharova-ID8bswh.html'
$ws.Cells.Item(100,7).Value = 'test field other'

# row 101
$ws.Cells.Item(101,1).Value = 'slando'
$ws.Cells.Item(101,2).Value = ', Приморский, Старосеннная площадь'
$ws.Cells.Item(101,3).Value = 'Выделенная коммуна на Старосеной площади'
$ws.Cells.Item(101,4).Value = '3 250 грн.'
$ws.Cells.Item(101,5).Value = 'Сегодня 13:36'
$ws.Cells.Item(101,6).Value = 'http://odessa.od.slando.ua/obyavlenie/vydelennaya-kommuna-na-starosenoy-ploschadi-ID87LC5.html'
$ws.Cells.Item(101,7).Value = 'test field other'

# row 102
$ws.Cells.Item(102,1).Value = 'slando'
$ws.Cells.Item(102,2).Value = ', Приморский, приморский'
$ws.Cells.Item(102,3).Value = 'Cдам 1-ком.в центре'
$ws.Cells.Item(102,4).Value = '2 500 грн.'
$ws.Cells.Item(102,5).Value = 'Сегодня 13:35'
$ws.Cells.Item(102,6).Value = 'http://odessa.od.slando.ua/obyavlenie/cdam-1-kom-v-tsentre-ID6ZOQP.html'
$ws.Cells.Item(102,7).Value = 'test field other'

# row 103
$ws.Cells.Item(103,1).Value = 'slando'
$ws.Cells.Item(103,2).Value = ', Ильфа и петрова 8'
$ws.Cells.Item(103,3).Value = 'Сдам 1 квартиру Ильфа и Петрова'
$ws.Cells.Item(103,4).Value = '2 600 грн.'
$ws.Cells.Item(103,5).Value = 'Сегодня 13:34'
$ws.Cells.Item(103,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-1-kvartiru-ilfa-i-petrova-ID59JXN.html'
$ws.Cells.Item(103,7).Value = 'test field other'

# row 104
$ws.Cells.Item(104,1).Value = 'slando'
$ws.Cells.Item(104,2).Value = ', Киевский, Люстдорфская дорога'
$ws.Cells.Item(104,3).Value = 'Сдам 2-х комнатную на люстдорфской дороге'
$ws.Cells.Item(104,4).Value = '4 069 грн.'
$ws.Cells.Item(104,5).Value = 'Сегодня 13:34'
$ws.Cells.Item(104,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-h-komnatnuyu-na-lyustdorfskoy-doroge-ID7YkFx.html'
$ws.Cells.Item(104,7).Value = 'test field other'

# row 105
$ws.Cells.Item(105,1).Value = 'slando'
$ws.Cells.Item(105,2).Value = ', Приморский, Ул. Кузнечная'
$ws.Cells.Item(105,3).Value = 'Сдам 1 комнатную квартиру'
$ws.Cells.Item(105,4).Value = '3 200 грн.'
$ws.Cells.Item(105,5).Value = 'Сегодня 13:33'
$ws.Cells.Item(105,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-1-komnatnuyu-kvartiru-ID7cn1r.html'
$ws.Cells.Item(105,7).Value = 'test field other'

# row 106
$ws.Cells.Item(106,1).Value = 'slando'
$ws.Cells.Item(106,2).Value = ', Приморский, большая арнаутская 57'
$ws.Cells.Item(106,3).Value = 'сдам 2 квартиру Большая арнаутская'
$ws.Cells.Item(106,4).Value = '3 200 грн.'
$ws.Cells.Item(106,5).Value = 'Сегодня 13:31'
$ws.Cells.Item(106,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-kvartiru-bolshaya-arnautskaya-ID7oDLn.html'
$ws.Cells.Item(106,7).Value = 'test field other'

# row 107
$ws.Cells.Item(107,1).Value = 'slando'
$ws.Cells.Item(107,2).Value = ', Приморский, Кленовая'
$ws.Cells.Item(107,3).Value = 'Сдам 2-х комнатную квартиру на Кленовой'
$ws.Cells.Item(107,4).Value = '5 697 грн.'
$ws.Cells.Item(107,5).Value = 'Сегодня 13:31'
$ws.Cells.Item(107,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-h-komnatnuyu-kvartiru-na-klenovoy-ID87IDb.html'
$ws.Cells.Item(107,7).Value = 'test field other'

# row 108
$ws.Cells.Item(108,1).Value = 'slando'
$ws.Cells.Item(108,2).Value = ', Рaзумовскя'
$ws.Cells.Item(108,3).Value = 'Сдам 2х комн.квартиру,с рaздельными ходaми'
$ws.Cells.Item(108,4).Value = '3 662 грн.'
$ws.Cells.Item(108,5).Value = 'Сегодня 13:29'
$ws.Cells.Item(108,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2h-komn-kvartiru-s-razdelnymi-hodami-ID4Yrbd.html'
$ws.Cells.Item(108,7).Value = 'test field other'

# row 109
$ws.Cells.Item(109,1).Value = 'slando'
$ws.Cells.Item(109,2).Value = ', Приморский, приморский р-н'
$ws.Cells.Item(109,3).Value = 'Сдам 2-Х ком.квартиру,в центре города,'
$ws.Cells.Item(109,4).Value = '4 883 грн.'
$ws.Cells.Item(109,5).Value = 'Сегодня 13:29'
$ws.Cells.Item(109,6).Value = 'http://odessa.od.slando.ua/obyavlenie/sdam-2-h-kom-kvartiru-v-tsentre-goroda-ID8b32l.html'
$ws.Cells.Item(109,7).Value = 'test field other'

# row 110
$ws.Cells.Item(110,1).Value = 'slando'
$ws.Cells.Item(110,2).Value = ''
$ws.Cells.Item(110,3).Value = 'Срочно квартира Базарная|ЖК Капитан,2раздельные спальни,+кухня'
$ws.Cells.Item(110,4).Value = '5 697 грн.'
$ws.Cells.Item(110,5).Value = 'Сегодня 13:28'
$ws.Cells.Item(110,6).Value = 'http://odessa.od.slando.ua/obyavlenie/srochno-kvartira-bazarnaya-zhk-kapitan-2razdelnye-spalni-kuhnya-ID8849d.html'
$ws.Cells.Item(110,7).Value = 'test field other'

# Rows 84-85 are highlighted in bold in the source workbook.
$ws.Range("A84:G85").Font.Bold = $true

# Restore the on-screen selection that was active when the workbook was last saved.
$ws.Range("A84:XFD85").Select()

# Add the new (empty) "test" worksheet after "Rent Data" and make it the active tab,
# matching the workbook-level activeTab="1" / tabSelected state in the saved file.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "test"
$newSheet.Activate()
